$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new employee row (row 95) that was imported via the new
#     "Import employee" feature referenced in the commit message. -----------
$ws.Cells.Item(95, 1).Value = 189
$ws.Cells.Item(95, 2).Value = "Test Import"
$ws.Cells.Item(95, 3).Value = 44753
$ws.Cells.Item(95, 4).Value = "abc"
$ws.Cells.Item(95, 5).Value = "abc"
$ws.Cells.Item(95, 6).Value = 35800

# Copy the date-formatted style from the row above (row 94) onto the new
# date-of-joining / date-of-birth cells so we reuse the existing cellXf
# (numFmtId 14) instead of minting a brand new one.
$ws.Range("C94").Copy()
$ws.Range("C95").PasteSpecial(-4122)
$ws.Range("F94").Copy()
$ws.Range("F95").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the view state: scrolled further down the list, new active
#     selection on the freshly-added row. -----------------------------------
$ws.Range("C95").Select()
$excel.ActiveWindow.ScrollRow = 85
$excel.ActiveWindow.ScrollColumn = 1
